# Updates cryptos list: refresh price / volume(1h) figures and restore the
# correct row ordering for two swapped coin pairs (PancakeSwap/Hedera and
# MXToken/Aave/HuobiToken/ARBITRUM), per the GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )
    # Route through a Text number-format so values that look numeric
    # (e.g. "1.00", "227.57") are not silently coerced into numbers,
    # then clear the format again so no stray style survives on save.
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "34.393.49"
Set-TextValue "E2" "  -0.16%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.800.96"
Set-TextValue "E3" "  +0.44%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.44%  "

# Row 5 - BNB
Set-TextValue "D5" "227.57"
Set-TextValue "E5" "  +0.42%  "

# Row 6 - XRP
Set-TextValue "D6" "0.577"
Set-TextValue "E6" "  +3.47%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.42%  "

# Row 8 - Solana
Set-TextValue "D8" "34.86"
Set-TextValue "E8" "  +5.69%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0690"
Set-TextValue "E10" "  -0.64%  "

# Row 11 - (Price unchanged)
Set-TextValue "E11" "  +0.20%  "

# Row 12 - Litecoin-ish big coin
Set-TextValue "D12" "2.058.51"
Set-TextValue "E12" "  +0.37%  "

# Row 13
Set-TextValue "D13" "11.19"
Set-TextValue "E13" "  +0.45%  "

# Row 14
Set-TextValue "D14" "1.807.86"
Set-TextValue "E14" "  -0.89%  "

# Row 15
Set-TextValue "D15" "0.641"
Set-TextValue "E15" "  +0.82%  "

# Row 16
Set-TextValue "D16" "34.364.19"
Set-TextValue "E16" "  -0.07%  "

# Row 17
Set-TextValue "D17" "4.32"
Set-TextValue "E17" "  +1.01%  "

# Row 18
Set-TextValue "D18" "69.06"
Set-TextValue "E18" "  +0.33%  "

# Row 19
Set-TextValue "D19" "244.50"
Set-TextValue "E19" "  -1.40%  "

# Row 20
Set-TextValue "D20" "0.0₃0793"
Set-TextValue "E20" "  -0.87%  "

# Row 21
Set-TextValue "D21" "11.48"
Set-TextValue "E21" "  +0.91%  "

# Row 22 - (Price unchanged)
Set-TextValue "E22" "  +0.48%  "

# Row 23
Set-TextValue "D23" "4.14"
Set-TextValue "E23" "  -0.67%  "

# Row 24
Set-TextValue "D24" "170.81"
Set-TextValue "E24" "  +3.54%  "

# Row 25 - (Price unchanged)
Set-TextValue "E25" "  +1.14%  "

# Row 26
Set-TextValue "D26" "7.50"
Set-TextValue "E26" "  +3.26%  "

# Row 27
Set-TextValue "D27" "16.68"
Set-TextValue "E27" "  +0.85%  "

# Row 28 - (Price unchanged)
Set-TextValue "E28" "  +1.50%  "

# Row 29
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.28%  "

# Row 30
Set-TextValue "D30" "3.97"
Set-TextValue "E30" "  +1.27%  "

# Row 31 / Row 32 - PancakeSwap and Hedera swap ranking positions
Set-TextValue "B31" "Hedera"
Set-TextValue "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.0527"
Set-TextValue "E31" "  +1.13%  "

Set-TextValue "B32" "PancakeSwap"
Set-TextValue "C32" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "1.24"
Set-TextValue "E32" "  +0.72%  "

# Row 33
Set-TextValue "D33" "3.80"
Set-TextValue "E33" "  -0.16%  "

# Row 34 - (Price unchanged)
Set-TextValue "E34" "  -0.08%  "

# Row 35
Set-TextValue "D35" "1.401.57"
Set-TextValue "E35" "  -1.21%  "

# Row 36
Set-TextValue "D36" "2.56"
Set-TextValue "E36" "  -1.38%  "

# Row 37
Set-TextValue "D37" "0.672"
Set-TextValue "E37" "  -0.13%  "

# Row 38 - (Price unchanged)
Set-TextValue "E38" "  +0.38%  "

# Row 39
Set-TextValue "D39" "0.0189"
Set-TextValue "E39" "  -1.78%  "

# Row 40 / Row 41 - MXToken and Aave swap ranking positions
Set-TextValue "B40" "Aave"
Set-TextValue "C40" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D40" "82.65"
Set-TextValue "E40" "  -2.72%  "

Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.84"
Set-TextValue "E41" "  +3.66%  "

# Row 42 / Row 43 - HuobiToken and ARBITRUM swap ranking positions
Set-TextValue "B42" "ARBITRUM"
Set-TextValue "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "0.948"
Set-TextValue "E42" "  +1.12%  "

Set-TextValue "B43" "HuobiToken"
Set-TextValue "C43" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D43" "2.41"
Set-TextValue "E43" "  +0.44%  "

# Row 44
Set-TextValue "D44" "13.72"
Set-TextValue "E44" "  +1.67%  "

# Row 45 - (Price unchanged)
Set-TextValue "E45" "  +2.75%  "

# Row 46
Set-TextValue "D46" "0.0512"
Set-TextValue "E46" "  -2.00%  "

# Row 47
Set-TextValue "D47" "5.99"
Set-TextValue "E47" "  -1.08%  "

# Row 48
Set-TextValue "D48" "1.960.77"
Set-TextValue "E48" "  +0.52%  "

# Row 49
Set-TextValue "D49" "104.74"
Set-TextValue "E49" "  -0.74%  "

# Row 50
Set-TextValue "D50" "1.00"
Set-TextValue "E50" "  +0.49%  "

# Row 51
Set-TextValue "D51" "0.0₆0129"
Set-TextValue "E51" "  +0.73%  "
